# Adds a new "290 packet run" data block (Acc/Loss/Time, columns I:K) to the
# "CNN 2D" worksheet, mirroring the 3rd group of columns already present on
# the "CNN 1D" sheet (header label + per-row Acc/Loss/Time values).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CNN 2D")

# --- headers -----------------------------------------------------------
# Row 1: group label over the new block (shares the existing shared-string
# "270 Packets" entry, same as sheet "CNN 1D" column I:K uses).
$ws.Range("J1").Value = "270 Packets"

# Row 2: Acc / Loss / Time column headers for the new block.
$ws.Range("I2").Value = "Acc"
$ws.Range("J2").Value = "Loss"
$ws.Range("K2").Value = "Time"

# --- data rows 3-51 (Acc, Loss, Time) -----------------------------------
$rows = @(
    @(3, 87.523418664932194, 0.27601916275212202, 86.042796850204397),
    @(4, 83.631747961044297, 0.28836688755879702, 82.178172588348303),
    @(5, 88.061755895614596, 0.25757784132007699, 82.852301836013794),
    @(6, 87.715512514114295, 0.27564888606440002, 78.845731258392306),
    @(7, 87.637251615524207, 0.27157388276724398, 84.234914541244507),
    @(8, 87.822228670120197, 0.27510275845858601, 83.965537786483694),
    @(9, 83.297365903854299, 0.32062178085626702, 83.761220932006793),
    @(10, 86.024618148803697, 0.298445907240941, 82.680081844329806),
    @(11, 88.358193635940495, 0.26434668743870698, 78.730302333831702),
    @(12, 87.215119600296006, 0.27926671390230401, 82.667416334152193),
    @(13, 87.573218345642005, 0.278621514167321, 83.788228511810303),
    @(14, 85.076010227203298, 0.28687788840489897, 78.902676582336397),
    @(15, 87.558990716934204, 0.265405901820637, 82.747016906738196),
    @(16, 83.491832017898503, 0.287447618680373, 83.847701549530001),
    @(17, 87.990605831146198, 0.28632474602824298, 78.339874029159503),
    @(18, 87.060970067977905, 0.27733732076754403, 78.372923135757404),
    @(19, 83.261793851852403, 0.30469947454019902, 82.768985271453801),
    @(20, 87.312352657318101, 0.31860044258782599, 78.103746414184499),
    @(21, 86.420661211013794, 0.29146146650996102, 83.461101055145207),
    @(22, 87.392985820770207, 0.28529822053480702, 84.087702035903902),
    @(23, 87.480729818344102, 0.28788178539928799, 82.849156379699707),
    @(24, 88.154244422912598, 0.26996309225579901, 84.929803609847994),
    @(25, 87.480729818344102, 0.26997012312736401, 78.031920909881507),
    @(26, 88.042783737182603, 0.26982108908240598, 83.901474714279104),
    @(27, 87.834089994430499, 0.27539024574641702, 83.4094624519348),
    @(28, 82.941639423370304, 2.6234215325207999, 81.862338066101003),
    @(29, 82.941639423370304, 2.6234215325207999, 82.8009259700775),
    @(30, 82.941639423370304, 2.6234215325207999, 82.688198804855304),
    @(31, 87.843573093414307, 0.268820609204653, 83.282687425613403),
    @(32, 82.941639423370304, 2.6234215325207999, 79.176611661910997),
    @(33, 87.094175815582204, 0.26979230514485403, 79.665502071380601),
    @(34, 86.420661211013794, 0.31910174135417702, 78.725725173950195),
    @(35, 87.7582013607025, 0.27807993861128399, 82.961811542510901),
    @(36, 83.302110433578406, 0.29837795085629798, 80.4127902984619),
    @(37, 87.570852041244507, 0.28663002960216899, 82.908334255218506),
    @(38, 87.103658914565997, 0.29996143048544099, 78.399708509445105),
    @(39, 87.6182794570922, 0.31409958712108998, 83.599899291992102),
    @(40, 87.134486436843801, 0.29030622976850001, 82.726615190505896),
    @(41, 87.253063917160006, 0.29968842456796302, 82.879879236221299),
    @(42, 87.687051296234102, 0.27863506997420601, 78.569346189498901),
    @(43, 82.941639423370304, 0.33626153936570802, 79.537481307983398),
    @(44, 88.0688667297363, 0.27671523363010903, 84.583001375198293),
    @(45, 87.748712301254201, 0.26505278520989101, 80.653891801834106),
    @(46, 87.646740674972506, 0.27481404390639003, 83.2636559009552),
    @(47, 87.518674135208101, 0.27853036461736802, 83.485069274902301),
    @(48, 87.736856937408405, 0.27262644007157699, 83.912900447845402),
    @(49, 87.094175815582204, 0.29761943380857298, 79.744165182113605),
    @(50, 87.900489568710299, 0.28061396861286497, 78.309705734252901),
    @(51, 82.941639423370304, 0.46228909113169803, 82.619557142257605)
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 5).Value = $row[1]   # column E - Acc
    $ws.Cells.Item($r, 6).Value = $row[2]   # column F - Loss
    $ws.Cells.Item($r, 7).Value = $row[3]   # column G - Time
}

# --- selection moved by the author after adding the data ---------------
$ws.Activate() | Out-Null
$ws.Range("M5").Select() | Out-Null
